$wb = $excel.ActiveWorkbook

# --- zh-cn sheet: bb251d22-... row (row 4) Correspond Handoff/Handback Datetime ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E4").Value = "2016-03-19 14:20:25"
$wsZh.Range("E5").Value = "2016-03-19 14:20:25"
$wsZh.Range("H4").Value = "2016-03-19 14:20:44"
$wsZh.Range("H5").Value = "2016-03-19 14:20:44"

# --- de-de sheet: bb251d22-... row (row 4) Correspond Handoff/Handback Datetime ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E4").Value = "2016-03-19 14:20:28"
$wsDe.Range("E5").Value = "2016-03-19 14:20:28"
$wsDe.Range("H4").Value = "2016-03-19 14:20:49"
$wsDe.Range("H5").Value = "2016-03-19 14:20:49"
